# feat: add single and multi corrector
#
# The original column A ("NO" / row index numbers) is removed and all
# remaining columns (代号/项目/结果/参考值/单位 ...) shift one column to the
# left. A handful of values in the remaining data are also corrected
# (unit spelling fixed to "μmol/L", a gamma prefix added to GGT's
# Chinese name, and a trailing zero trimmed from the glucose result).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column A entirely; this shifts B:F left into A:E and updates
# the used range / dimension accordingly.
$ws.Columns.Item(1).Delete()

# Correct a few remaining values after the shift.
$ws.Range("E2").Value = "μmol/L"
$ws.Range("E3").Value = "μmol/L"
$ws.Range("E4").Value = "μmol/L"
$ws.Range("B7").Value = "γ谷氨酰转肽酶"
$ws.Range("E13").Value = "μmol/L"

# Keep the trimmed glucose result as text (matching the rest of the
# sheet, which stores every value as text) instead of letting Excel
# auto-convert "4.2" into a numeric value.
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "4.2"
